$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Propagate the existing row 147 formatting (styles only) to the new
# rows 148-150 before any values are written, so cell styles (s="1"
# on column A, s="2" on column E) match the rest of the sheet.
# Row 148 needs the full A:AC range (it has H/I/J/AB/AC populated),
# while rows 149 and 150 only need A:G and K:AA (no H/I/J/AB/AC).
# ------------------------------------------------------------------
$ws.Range("A147:AC147").Copy()
$ws.Range("A148:AC148").PasteSpecial(-4122)

$ws.Range("A147:G147").Copy()
$ws.Range("A149:G149").PasteSpecial(-4122)
$ws.Range("A150:G150").PasteSpecial(-4122)

$ws.Range("K147:AA147").Copy()
$ws.Range("K149:AA149").PasteSpecial(-4122)
$ws.Range("K150:AA150").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Row 148 (new match: id 146)
# ------------------------------------------------------------------
$ws.Range("A148").Value = 146
$ws.Range("B148").Value = 7751751
$ws.Range("C148").Value = "India Super League"
$ws.Range("D148").Value = "India Super League"
$ws.Range("E148").Value = 45346.45833333334
$ws.Range("F148").Value = "Bengaluru"
$ws.Range("G148").Value = "Hyderabad FC"
$ws.Range("H148").Value = 2
$ws.Range("I148").Value = 1
$ws.Range("J148").Value = "H"
$ws.Range("K148").Value = 1.2
$ws.Range("L148").Value = 5.5
$ws.Range("M148").Value = 10
$ws.Range("N148").Value = 1.2
$ws.Range("O148").Value = 5.25
$ws.Range("P148").Value = 12
$ws.Range("Q148").Value = -1.75
$ws.Range("R148").Value = 1.875
$ws.Range("S148").Value = 1.975
$ws.Range("T148").Value = 3
$ws.Range("U148").Value = 2.025
$ws.Range("V148").Value = 1.825
$ws.Range("W148").Value = 0.2
$ws.Range("X148").Value = -1
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = -1
$ws.Range("AA148").Value = 0.9750000000000001
$ws.Range("AB148").Value = 0
$ws.Range("AC148").Value = -0

# ------------------------------------------------------------------
# Row 149 (moved original row 147 match: id 147, odds U/V refreshed)
# ------------------------------------------------------------------
$ws.Range("A149").Value = 147
$ws.Range("B149").Value = 7749466
$ws.Range("C149").Value = "India Super League"
$ws.Range("D149").Value = "India Super League"
$ws.Range("E149").Value = 45347.45833333334
$ws.Range("F149").Value = "Kerala Blasters"
$ws.Range("G149").Value = "FC Goa"
$ws.Range("K149").Value = 3.3
$ws.Range("L149").Value = 3.4
$ws.Range("M149").Value = 1.95
$ws.Range("N149").Value = 3.5
$ws.Range("O149").Value = 3.4
$ws.Range("P149").Value = 1.909
$ws.Range("Q149").Value = 0.5
$ws.Range("R149").Value = 1.8
$ws.Range("S149").Value = 2
$ws.Range("T149").Value = 2.5
$ws.Range("U149").Value = 1.875
$ws.Range("V149").Value = 1.925
$ws.Range("W149").Value = 0
$ws.Range("X149").Value = 0
$ws.Range("Y149").Value = 0
$ws.Range("Z149").Value = 0
$ws.Range("AA149").Value = 0

# ------------------------------------------------------------------
# Row 150 (new match: id 148)
# ------------------------------------------------------------------
$ws.Range("A150").Value = 148
$ws.Range("B150").Value = 7751752
$ws.Range("C150").Value = "India Super League"
$ws.Range("D150").Value = "India Super League"
$ws.Range("E150").Value = 45348.45833333334
$ws.Range("F150").Value = "East Bengal Club"
$ws.Range("G150").Value = "Chennaiyin FC"
$ws.Range("K150").Value = 2.3
$ws.Range("L150").Value = 3.2
$ws.Range("M150").Value = 2.8
$ws.Range("N150").Value = 2.3
$ws.Range("O150").Value = 3.2
$ws.Range("P150").Value = 2.8
$ws.Range("Q150").Value = -0.25
$ws.Range("R150").Value = 2.05
$ws.Range("S150").Value = 1.75
$ws.Range("T150").Value = 2.5
$ws.Range("U150").Value = 1.875
$ws.Range("V150").Value = 1.925
$ws.Range("W150").Value = 0
$ws.Range("X150").Value = 0
$ws.Range("Y150").Value = 0
$ws.Range("Z150").Value = 0
$ws.Range("AA150").Value = 0

# ------------------------------------------------------------------
# Row 147 (overwritten in place with a different match: id 145)
# ------------------------------------------------------------------
$ws.Range("B147").Value = 7873049
$ws.Range("E147").Value = 45346.35416666666
$ws.Range("F147").Value = "Odisha FC"
$ws.Range("G147").Value = "Mohun Bagan SG"
$ws.Range("H147").Value = 0
$ws.Range("I147").Value = 0
$ws.Range("J147").Value = "D"
$ws.Range("K147").Value = 2.25
$ws.Range("L147").Value = 3.5
$ws.Range("M147").Value = 2.625
$ws.Range("N147").Value = 2.375
$ws.Range("O147").Value = 3.1
$ws.Range("P147").Value = 2.7
$ws.Range("Q147").Value = 0
$ws.Range("R147").Value = 1.775
$ws.Range("S147").Value = 2.1
$ws.Range("T147").Value = 2.5
$ws.Range("U147").Value = 1.9
$ws.Range("V147").Value = 1.95
$ws.Range("W147").Value = -1
$ws.Range("X147").Value = 2.1
$ws.Range("Y147").Value = -1
$ws.Range("Z147").Value = 0
$ws.Range("AA147").Value = -0
$ws.Range("AB147").Value = -1
$ws.Range("AC147").Value = 0.95
